$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column B ("Taxonsorteringsordning") gets bumped by +4 for every data row
# from row 42 through row 80. Rows 65/66 and 78/79 additionally have their
# entire record (all other columns) swapped between the pair, because the
# update reordered which species occupies which row; those four rows are
# therefore handled explicitly below instead of via the simple +4 rule.
# ---------------------------------------------------------------------------

$simpleRows = @(42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,67,68,69,70,71,72,73,74,75,76,77,80)

$newB = @{
    42 = 79243
    43 = 79243
    44 = 79243
    45 = 78255
    46 = 79243
    47 = 79243
    48 = 91804
    49 = 79243
    50 = 79243
    51 = 91804
    52 = 79243
    53 = 78255
    54 = 79243
    55 = 83228
    56 = 79243
    57 = 79243
    58 = 91828
    59 = 79243
    60 = 79243
    61 = 79243
    62 = 79243
    63 = 79243
    64 = 78255
    67 = 79243
    68 = 79243
    69 = 83223
    70 = 78255
    71 = 83223
    72 = 79243
    73 = 83223
    74 = 79499
    75 = 83223
    76 = 89193
    77 = 79243
    80 = 79243
}

foreach ($r in $simpleRows) {
    $ws.Range("B$r").Value = $newB[$r]
}

# ---------------------------------------------------------------------------
# Row 65 <-> Row 66 swap (with the new +4 taxon sort order baked in)
# ---------------------------------------------------------------------------

$ws.Range("A65").Value = 130837548
$ws.Range("B65").Value = 83223
$ws.Range("D65").Value = "NT"
$ws.Range("E65").Value = 6440
$ws.Range("F65").Value = "Vitgrynig nållav"
$ws.Range("G65").Value = "Chaenotheca subroscida"
$ws.Range("H65").Value = "(Eitner) Zahlbr."
$ws.Range("Q65").Value = 445740
$ws.Range("R65").Value = 7026322
$ws.Range("S65").Value = 8
$ws.Range("Z65").Value = "10:53"
$ws.Range("AB65").Value = "10:53"
$ws.Range("AC65").Value = "På bark på stam av levande gammal gran i gammal granskog"

$ws.Range("A66").Value = 130839350
$ws.Range("B66").Value = 79243
$ws.Range("D66").Value = "NT"
$ws.Range("E66").Value = 6425
$ws.Range("F66").Value = "Garnlav"
$ws.Range("G66").Value = "Alectoria sarmentosa"
$ws.Range("H66").Value = "(Ach.) Ach."
$ws.Range("Q66").Value = 445790
$ws.Range("R66").Value = 7026340
$ws.Range("S66").Value = 6
$ws.Range("Z66").Value = "12:47"
$ws.Range("AB66").Value = "12:47"
$ws.Range("AC66").Value = "På gammal gran i gammal granskog"

# ---------------------------------------------------------------------------
# Row 78 <-> Row 79 swap (with the new taxon sort order baked in)
# ---------------------------------------------------------------------------

$ws.Range("A78").Value = 130839413
$ws.Range("B78").Value = 78255
$ws.Range("D78").Value = "NT"
$ws.Range("E78").Value = 228579
$ws.Range("F78").Value = "Liten svartspik"
$ws.Range("G78").Value = "Chaenothecopsis nana"
$ws.Range("H78").Value = "Tibell"
$ws.Range("Q78").Value = 445781
$ws.Range("R78").Value = 7026373
$ws.Range("S78").Value = 7
$ws.Range("Z78").Value = "12:54"
$ws.Range("AB78").Value = "12:54"
$ws.Range("AC78").Value = "På bark på stam av levande gammal gran i gles gammal granskog"

$ws.Range("A79").Value = 130837541
$ws.Range("B79").Value = 75221
$ws.Range("D79").Value = "LC"
$ws.Range("E79").Value = 6428
$ws.Range("F79").Value = "Rostfläck"
$ws.Range("G79").Value = "Arthonia vinosa"
$ws.Range("H79").Value = "Leight."
$ws.Range("Q79").Value = 445740
$ws.Range("R79").Value = 7026322
$ws.Range("S79").Value = 8
$ws.Range("Z79").Value = "10:52"
$ws.Range("AB79").Value = "10:52"
$ws.Range("AC79").Value = "På tunna kvistar vid basen på gammal levande gran"
